# Update jurnal & templates
# Adds one new journal entry (row 13) to the "SUMBER" table: title, date,
# and OFFLINE/ONLINE hyperlink cells, then leaves the selection where the
# author ended up after typing the row (one cell past the table, F13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: new journal entry -------------------------------------------------
$ws.Range("B13").Value = "PENERJEMAH BAHASA ALAMI DALAM BAHASA INDONESIA KE SOURCE CODE DALAM BAHASA PASCAL"
$ws.Range("C13").Value = 45273   # 12/13/2023

# OFFLINE (D13) -> local PDF copy, same convention as the other rows.
$ws.Hyperlinks.Add($ws.Range("D13"), "PENERJEMAH%20BAHASA%20ALAMI%20DALAM%20BAHASA%20INDONESIA%20KE%20SOURCE%20CODE%20DALAM%20BAHASA%20PASCAL.pdf", $null, $null, "GET")
# Re-apply the donor cell's number/border formatting so D13 reuses the
# existing "hyperlink + bordered" style instead of a brand-new one.
$ws.Range("D3").Copy()
$ws.Range("D13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ONLINE (E13) -> external source for the paper.
$ws.Hyperlinks.Add($ws.Range("E13"), "https://ejournal.stmik-time.ac.id/index.php/jurnalTIME/article/view/penerjemah-bahasa-alami-pascal", $null, $null, "GET")
$ws.Range("E3").Copy()
$ws.Range("E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# --- Final cursor position, matching the author's saved selection --------------
$ws.Range("F13").Select()
